# bug fix scan & open fail
# Fill in bug report row 12 (issue #11) on the buglist sheet and move the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 12: new bug entry -------------------------------------------------
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "부팅 후, 인증안됨."
$ws.Range("C12").Value = "항상"
$ws.Range("D12").Value = "높음"
$ws.Range("E12").Value = "검증"

# F12 / G12 are dates; copy the number format from the row above (style reuse)
# before writing the serial date values so the date format is preserved.
$ws.Range("F11").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("F12").Value = 42012
$ws.Range("G12").Value = 42014

$ws.Range("H12").Value = "#2 를 수정하며 생긴 사이드이펙트. Update할 내용이 없을 경우, scan을 시작하는 코드가 누락되었었습니다. (2015-01-10)"

# Row 12 grows to fit the wrapped remarks text
$ws.Rows.Item(12).RowHeight = 33

# --- Move active selection to F17 ------------------------------------------
$null = $ws.Range("F17").Select()
